# alcustoms.Excel Tables update (#33)
# Add a new "Tables4" worksheet (after "Tables3") that exercises the new
# greedycolumns / greedyrows behavior of excel.Tables.gettablesize, and make
# it the active sheet/tab - mirroring the author's manual Excel session.

$wb = $excel.ActiveWorkbook

# Insert the new worksheet at the end of the sheet list (after the current
# last sheet, "Tables3") and give it its final name.
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Tables4"

# --- "We Format This is Test 1" block (A1:C3) -----------------------------
$ws.Range("A1").Value = "We"
$ws.Range("C1").Value = "Format"
# --- "Our  Tables / Funny" greedy-column demo block (F1:G3) ---------------
$ws.Range("F1").Value = "Our  Tables"
$ws.Range("G1").Value = "Funny"
$ws.Range("A2").Value = "This"
$ws.Range("C2").Value = "is"
$ws.Range("A3").Value = "Test"
$ws.Range("C3").Value = 1
$ws.Range("F2").Value = "Test 2"
$ws.Range("G2").Value = "Should"
$ws.Range("F3").Value = "End"
$ws.Range("G3").Value = "Here"

# --- "We Have Greedy Rows" greedy-row demo block (K1:N7) -------------------
$ws.Range("K1").Value = "We"
$ws.Range("L1").Value = "Have"
$ws.Range("M1").Value = "Greedy"
$ws.Range("N1").Value = "Rows"

$ws.Range("K2").Value = "This"
$ws.Range("L2").Value = "row"
$ws.Range("M2").Value = "is "
$ws.Range("N2").Value = "Normal"

$ws.Range("K4").Value = "This"
$ws.Range("L4").Value = "row"
$ws.Range("M4").Value = "is "
$ws.Range("N4").Value = "Greedy"

$ws.Range("K7").Value = "This"
$ws.Range("L7").Value = "row"
$ws.Range("M7").Value = "is"
$ws.Range("N7").Value = "Super Greedy"

# Column F was widened (best-fit) for the "Our  Tables"/"Test 2" labels.
$ws.Columns.Item(6).ColumnWidth = 10

# Leave the selection where the author's session ended up, and make this
# new sheet the active tab (matches bookViews/activeTab + sheetView
# tabSelected moving from "Tables3" to "Tables4").
$ws.Range("L8").Select() | Out-Null
$ws.Activate()
